$d = $word.ActiveDocument

$d.Content.Find.Execute("526÷6=87, 4", $true, $false, $false, $false, $false, $true, 1, $false, "948÷6=158, 0", 2)
$d.Content.Find.Execute("484÷9=53, 7", $true, $false, $false, $false, $false, $true, 1, $false, "673÷6=112, 1", 2)
$d.Content.Find.Execute("657÷9=73, 0", $true, $false, $false, $false, $false, $true, 1, $false, "125÷3=41, 2", 2)
$d.Content.Find.Execute("624÷9=69, 3", $true, $false, $false, $false, $false, $true, 1, $false, "883÷6=147, 1", 2)
$d.Content.Find.Execute("585÷4=146, 1", $true, $false, $false, $false, $false, $true, 1, $false, "841÷4=210, 1", 2)
$d.Content.Find.Execute("379÷5=75, 4", $true, $false, $false, $false, $false, $true, 1, $false, "276÷5=55, 1", 2)
$d.Content.Find.Execute("508÷7=72, 4", $true, $false, $false, $false, $false, $true, 1, $false, "961÷9=106, 7", 2)
$d.Content.Find.Execute("613÷6=102, 1", $true, $false, $false, $false, $false, $true, 1, $false, "136÷7=19, 3", 2)
$d.Content.Find.Execute("142÷5=28, 2", $true, $false, $false, $false, $false, $true, 1, $false, "621÷4=155, 1", 2)
$d.Content.Find.Execute("969÷3=323, 0", $true, $false, $false, $false, $false, $true, 1, $false, "813÷7=116, 1", 2)
$d.Content.Find.Execute("241÷5=48, 1", $true, $false, $false, $false, $false, $true, 1, $false, "244÷5=48, 4", 2)
$d.Content.Find.Execute("674÷8=84, 2", $true, $false, $false, $false, $false, $true, 1, $false, "435÷7=62, 1", 2)
$d.Content.Find.Execute("231÷7=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "566÷3=188, 2", 2)
$d.Content.Find.Execute("156÷7=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "688÷8=86, 0", 2)
$d.Content.Find.Execute("163÷3=54, 1", $true, $false, $false, $false, $false, $true, 1, $false, "521÷6=86, 5", 2)
$d.Content.Find.Execute("566÷6=94, 2", $true, $false, $false, $false, $false, $true, 1, $false, "765÷4=191, 1", 2)
$d.Content.Find.Execute("460÷8=57, 4", $true, $false, $false, $false, $false, $true, 1, $false, "513÷4=128, 1", 2)
$d.Content.Find.Execute("961÷4=240, 1", $true, $false, $false, $false, $false, $true, 1, $false, "143÷3=47, 2", 2)
$d.Content.Find.Execute("214÷8=26, 6", $true, $false, $false, $false, $false, $true, 1, $false, "910÷8=113, 6", 2)
$d.Content.Find.Execute("137÷9=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "919÷7=131, 2", 2)
$d.Content.Find.Execute("332÷5=66, 2", $true, $false, $false, $false, $false, $true, 1, $false, "108÷8=13, 4", 2)
$d.Content.Find.Execute("141÷9=15, 6", $true, $false, $false, $false, $false, $true, 1, $false, "751÷6=125, 1", 2)
$d.Content.Find.Execute("986÷3=328, 2", $true, $false, $false, $false, $false, $true, 1, $false, "235÷3=78, 1", 2)
$d.Content.Find.Execute("794÷2=397, 0", $true, $false, $false, $false, $false, $true, 1, $false, "981÷2=490, 1", 2)
$d.Content.Find.Execute("513÷7=73, 2", $true, $false, $false, $false, $false, $true, 1, $false, "602÷7=86, 0", 2)
